$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(2, 1).Range.Text = "Max Verstappen"
$t.Cell(2, 4).Range.Text = "1"
$t.Cell(2, 5).Range.Text = "-2"
$t.Cell(2, 6).Range.Text = "235.606"
$t.Cell(2, 7).Range.Text = "0"
$t.Cell(3, 1).Range.Text = "Lando Norris"
$t.Cell(3, 4).Range.Text = "2"
$t.Cell(3, 6).Range.Text = "234.256"
$t.Cell(3, 7).Range.Text = "-1.350"
$t.Cell(4, 5).Range.Text = "-5"
$t.Cell(4, 6).Range.Text = "233.961"
$t.Cell(4, 7).Range.Text = "-1.644"
$t.Cell(5, 5).Range.Text = "-7"
$t.Cell(5, 6).Range.Text = "230.525"
$t.Cell(5, 7).Range.Text = "-5.081"
$t.Cell(6, 5).Range.Text = "-2"
$t.Cell(6, 6).Range.Text = "228.629"
$t.Cell(6, 7).Range.Text = "-6.976"
$t.Cell(7, 5).Range.Text = "-4"
$t.Cell(7, 6).Range.Text = "237.660"
$t.Cell(7, 7).Range.Text = "2.055"
$t.Cell(8, 5).Range.Text = "-4"
$t.Cell(8, 6).Range.Text = "235.065"
$t.Cell(8, 7).Range.Text = "-0.541"
$t.Cell(9, 1).Range.Text = "Sergio Perez"
$t.Cell(9, 4).Range.Text = "8"
$t.Cell(9, 5).Range.Text = "-3"
$t.Cell(9, 6).Range.Text = "224.046"
$t.Cell(9, 7).Range.Text = "-11.560"
$t.Cell(10, 1).Range.Text = "Pierre Gasly"
$t.Cell(10, 4).Range.Text = "9"
$t.Cell(10, 5).Range.Text = "-7"
$t.Cell(10, 6).Range.Text = "235.294"
$t.Cell(10, 7).Range.Text = "-0.311"
$t.Cell(11, 5).Range.Text = "-7"
$t.Cell(11, 7).Range.Text = "-2.065"
$t.Cell(12, 6).Range.Text = "243.132"
$t.Cell(12, 7).Range.Text = "7.526"
$t.Cell(13, 5).Range.Text = "-6"
$t.Cell(13, 6).Range.Text = "309.969"
$t.Cell(13, 7).Range.Text = "-5.411"
$t.Cell(14, 5).Range.Text = "-6"
$t.Cell(14, 6).Range.Text = "310.341"
$t.Cell(14, 7).Range.Text = "-5.039"
$t.Cell(15, 5).Range.Text = "-8"
$t.Cell(15, 6).Range.Text = "286.892"
$t.Cell(15, 7).Range.Text = "-28.487"
$t.Cell(16, 5).Range.Text = "-4"
$t.Cell(16, 6).Range.Text = "310.109"
$t.Cell(16, 7).Range.Text = "-5.270"
$t.Cell(17, 1).Range.Text = "Alexander Albon"
$t.Cell(17, 4).Range.Text = "18"
$t.Cell(17, 5).Range.Text = "-5"
$t.Cell(17, 6).Range.Text = "301.461"
$t.Cell(17, 7).Range.Text = "-13.918"
$t.Cell(18, 5).Range.Text = "-4"
$t.Cell(18, 6).Range.Text = "307.505"
$t.Cell(18, 7).Range.Text = "-7.874"
$t.Cell(19, 1).Range.Text = "Guanyu Zhou"
$t.Cell(19, 4).Range.Text = "13"
$t.Cell(19, 5).Range.Text = "-8"
$t.Cell(19, 6).Range.Text = "359.138"
$t.Cell(19, 7).Range.Text = "43.758"
$t.Cell(20, 5).Range.Text = "-8"
$t.Cell(20, 6).Range.Text = "320.763"
$t.Cell(20, 7).Range.Text = "5.383"
$t.Cell(21, 5).Range.Text = "-8"
$t.Cell(21, 6).Range.Text = "394.399"
$t.Cell(21, 7).Range.Text = "-0.447"
